$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tactics")

# New x-mitre-tactic STIX IDs (ATT&CK v14.1 / ATLAS data 4.5.1 refresh), in row order (rows 2-15)
$newIds = @(
    "x-mitre-tactic--afe18db5-cd3b-404c-8993-3aaf488c5f62",
    "x-mitre-tactic--f3deea7c-6d5a-49e7-ae3a-d70afc7889c9",
    "x-mitre-tactic--9633dd67-f5e8-4a19-88aa-c9b6fcf9596c",
    "x-mitre-tactic--c8c04179-cc2c-4a1d-b392-e6d4c42c3493",
    "x-mitre-tactic--d9deeb4d-c193-469d-a1c8-52e6d145da0c",
    "x-mitre-tactic--2538ccce-34dd-4470-8c3f-747be088b4e6",
    "x-mitre-tactic--f7d86c97-d8a6-4b81-a4a3-b49191faeced",
    "x-mitre-tactic--876860aa-e340-435c-824b-b5dd011c4c95",
    "x-mitre-tactic--c05767a2-5965-4618-b66e-6a7fffed256a",
    "x-mitre-tactic--813efb61-5bed-439d-a264-1b4e0ab89c37",
    "x-mitre-tactic--8bb9a040-d075-4566-966d-4a6a588f4fad",
    "x-mitre-tactic--39d150db-6130-4b2e-81e5-754b344a6a00",
    "x-mitre-tactic--d4d10ed9-f6d7-4774-a9f0-68bf4970c5bd",
    "x-mitre-tactic--4cd8009a-c196-4c3d-8eae-87afa661ed8e"
)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newIds[$i]
}

# Update the "created"/"last modified" date (column F and G) for all data rows
$newDate = "12 January 2024"
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 6).Value = $newDate
    $ws.Cells.Item($row, 7).Value = $newDate
}
